$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: 2025-04-03, Scrabble -> 2025-05-29, Catan
$ws.Range("B22").Value = "'2025-05-29"
$ws.Range("C22").Value = "Catan"
$ws.Range("D22").Value = 34
$ws.Range("E22").Value = 44
$ws.Range("F22").Value = 123
$ws.Range("G22").Value = 47
$ws.Range("H22").Value = 57

# Row 23: 2025-05-20, Ticket to Ride -> 2025-06-03, Catan
$ws.Range("B23").Value = "'2025-06-03"
$ws.Range("C23").Value = "Catan"
$ws.Range("D23").Value = 31
$ws.Range("E23").Value = 444
$ws.Range("F23").Value = 152
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 42
$ws.Range("J23").Value = "morning tournament"

# Row 24: 2025-05-01 -> 2025-06-01 (name_of_game stays Scrabble)
$ws.Range("B24").Value = "'2025-06-01"
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = 56
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 199
$ws.Range("H24").Value = 336

# Row 25: 2025-05-01, Ticket to Ride -> 2025-06-19, Chess
$ws.Range("B25").Value = "'2025-06-19"
$ws.Range("C25").Value = "Chess"
$ws.Range("D25").Value = 47
$ws.Range("E25").Value = 391
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 7
$ws.Range("J25").Value = "meet 30 minutes before start"

# Row 26: 2025-05-14, Monopoly -> 2025-05-04, Scrabble
$ws.Range("B26").Value = "'2025-05-04"
$ws.Range("C26").Value = "Scrabble"
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = 21
$ws.Range("F26").Value = 702
$ws.Range("G26").Value = 195
$ws.Range("H26").Value = 21
$ws.Range("J26").Value = "evening tournament"
